# docs/version_scheduling/0.0.11.xlsx - "split out gamemode specific code"
#
# The "Fix relative velocity being added multiple times" task (row 6) and the
# "Alt+tab crash in dedicated fullscreen" task (row 12) are marked as
# Completed ("Yes") in column C, highlighted with the accent6 theme fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark both bugfix rows as completed.
$ws.Range("C6").Value = "Yes"
$ws.Range("C12").Value = "Yes"

# Highlight C6 with the green (accent6/theme 9) fill used for "done" items.
$ws.Range("C6").Interior.ThemeColor = 10

# Copy that formatting onto C12 so both cells share the exact same style.
$ws.Range("C6").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection where the author ended up after the edit.
$ws.Range("A16").Select() | Out-Null
